# Refresh the cryptocurrency price/volume snapshot (and fix the ImmutableX /
# ARBITRUM row ordering swap) to match the latest scrape.
#
# D (Price) and E (Volume(1h)) hold plain text in the source workbook (e.g.
# "1.000", "0.9999", "  -1.04%  ") - NOT numbers. Assigning such strings
# straight to Range.Value would let Excel's input-parsing "helpfully" treat
# them as numbers/dates, silently dropping significant trailing zeros
# (e.g. "1.000" -> 1) and changing the stored cell type. To guarantee every
# written value lands back in the sheet as the exact literal text (with no
# left-over formulas and no incidental style/number-format changes), each
# value is first pushed in through a `="..."` text formula - which always
# evaluates to a plain string - and the whole touched block is then
# converted from formulas to static values via Copy + PasteSpecial(values).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Formula = '="26.853.24"'
$ws.Range("E2").Formula = '="  -1.07%  "'
$ws.Range("D3").Formula = '="1.804.91"'
$ws.Range("E3").Formula = '="  -0.90%  "'
$ws.Range("D4").Formula = '="1.000"'
$ws.Range("E4").Formula = '="  -0.44%  "'
$ws.Range("D5").Formula = '="310.01"'
$ws.Range("E5").Formula = '="  -0.76%  "'
$ws.Range("D6").Formula = '="0.9999"'
$ws.Range("E6").Formula = '="  -0.45%  "'
$ws.Range("D7").Formula = '="0.4449"'
$ws.Range("E7").Formula = '="  +5.15%  "'
$ws.Range("D8").Formula = '="0.3668"'
$ws.Range("E8").Formula = '="  -0.40%  "'
$ws.Range("D9").Formula = '="0.07335"'
$ws.Range("E9").Formula = '="  +1.41%  "'
$ws.Range("D10").Formula = '="0.8562"'
$ws.Range("E10").Formula = '="  +0.20%  "'
$ws.Range("D11").Formula = '="20.61"'
$ws.Range("E11").Formula = '="  -1.48%  "'
$ws.Range("D12").Formula = '="1.935.77"'
$ws.Range("E12").Formula = '="  +6.16%  "'
$ws.Range("D13").Formula = '="6.596"'
$ws.Range("E13").Formula = '="  -1.38%  "'
$ws.Range("D14").Formula = '="92.65"'
$ws.Range("E14").Formula = '="  +3.68%  "'
$ws.Range("D15").Formula = '="0.07069"'
$ws.Range("E15").Formula = '="  -0.25%  "'
$ws.Range("D16").Formula = '="5.290"'
$ws.Range("E16").Formula = '="  +0.02%  "'
$ws.Range("E17").Formula = '="  -0.47%  "'
$ws.Range("D18").Formula = '="0.000008718"'
$ws.Range("E18").Formula = '="  -1.35%  "'
$ws.Range("E19").Formula = '="  -0.43%  "'
$ws.Range("D20").Formula = '="14.85"'
$ws.Range("E20").Formula = '="  -1.02%  "'
$ws.Range("D21").Formula = '="26.888.18"'
$ws.Range("E21").Formula = '="  -1.20%  "'
$ws.Range("D22").Formula = '="5.144"'
$ws.Range("E22").Formula = '="  +0.51%  "'
$ws.Range("E23").Formula = '="  -0.49%  "'
$ws.Range("D24").Formula = '="1.992"'
$ws.Range("E24").Formula = '="  +0.42%  "'
$ws.Range("D25").Formula = '="151.85"'
$ws.Range("E25").Formula = '="  -0.32%  "'
$ws.Range("D26").Formula = '="18.41"'
$ws.Range("E26").Formula = '="  +0.18%  "'
$ws.Range("D27").Formula = '="2.169"'
$ws.Range("E27").Formula = '="  -0.82%  "'
$ws.Range("D28").Formula = '="5.194"'
$ws.Range("E28").Formula = '="  -0.43%  "'
$ws.Range("D29").Formula = '="116.72"'
$ws.Range("E29").Formula = '="  +0.55%  "'
$ws.Range("D30").Formula = '="0.08821"'
$ws.Range("E30").Formula = '="  -0.19%  "'
$ws.Range("B31").Formula = '="ARBITRUM"'
$ws.Range("C31").Formula = '="https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"'
$ws.Range("D31").Formula = '="1.170"'
$ws.Range("E31").Formula = '="  -1.26%  "'
$ws.Range("B32").Formula = '="ImmutableX"'
$ws.Range("C32").Formula = '="https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"'
$ws.Range("D32").Formula = '="0.7459"'
$ws.Range("E32").Formula = '="  -0.06%  "'
$ws.Range("D33").Formula = '="2.930"'
$ws.Range("E33").Formula = '="  -1.58%  "'
$ws.Range("D34").Formula = '="4.454"'
$ws.Range("E34").Formula = '="  +0.45%  "'
$ws.Range("D35").Formula = '="0.9994"'
$ws.Range("E35").Formula = '="  -0.57%  "'
$ws.Range("E36").Formula = '="  -2.43%  "'
$ws.Range("D37").Formula = '="0.01963"'
$ws.Range("E37").Formula = '="  +0.00%  "'
$ws.Range("D38").Formula = '="0.05186"'
$ws.Range("E38").Formula = '="  -0.84%  "'
$ws.Range("D39").Formula = '="0.5311"'
$ws.Range("E39").Formula = '="  +5.73%  "'
$ws.Range("E40").Formula = '="  -0.54%  "'
$ws.Range("D41").Formula = '="7.015"'
$ws.Range("E41").Formula = '="  -3.73%  "'
$ws.Range("D42").Formula = '="0.1681"'
$ws.Range("E42").Formula = '="  -0.68%  "'
$ws.Range("D43").Formula = '="0.5124"'
$ws.Range("E43").Formula = '="  +8.26%  "'
$ws.Range("D44").Formula = '="8.454"'
$ws.Range("E44").Formula = '="  -2.14%  "'
$ws.Range("D45").Formula = '="10.51"'
$ws.Range("E45").Formula = '="  -0.47%  "'
$ws.Range("D46").Formula = '="1.970"'
$ws.Range("E46").Formula = '="  +5.96%  "'
$ws.Range("D47").Formula = '="105.44"'
$ws.Range("E47").Formula = '="  -0.73%  "'
$ws.Range("D48").Formula = '="0.9989"'
$ws.Range("E48").Formula = '="  -0.56%  "'
$ws.Range("D49").Formula = '="1.662"'
$ws.Range("E49").Formula = '="  +0.18%  "'
$ws.Range("E50").Formula = '="  -1.13%  "'
$ws.Range("D51").Formula = '="0.9161"'
$ws.Range("E51").Formula = '="  +0.57%  "'

# Collapse the temporary "="..."" formulas down to plain static text values
# (B2:E51 is a superset covering every cell touched above, including the
# B31:C32 coin-name/link swap) so the saved file has no formulas and no
# stray formatting left behind.
$touched = $ws.Range("B2:E51")
$touched.Copy($touched)
$touched.PasteSpecial(-4163)
